$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Marking" row (row 11): Right marks and Wrong (negative) marks
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Update "Total" row (row 12): Right total, Wrong total, and Max summary text
$ws.Range("B12").Value = 50
$ws.Range("C12").Value = -0
$ws.Range("E12").Value = "50.0/140"
